$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Step 1: Insert a new "2022-Q4" row at the top of the "总计" summary
#         sheet's data (row 2), pushing the existing rows down by one.
# -----------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The insert copies the formatting of the row above (the bold/bordered
# header) into the new row; restore the correct per-column formatting by
# copying it from the row directly below (which already has the right,
# shifted-down formatting of a normal data row).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B3:D3").Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.62

# -----------------------------------------------------------------------
# Step 2: Add the new "2022-Q4" fund-holdings sheet, positioned right
#         after "总计" and before "2022-Q3". Duplicate the "2022-Q3"
#         sheet (so the header row / column styling matches the other
#         quarterly sheets exactly) and then overwrite its data.
# -----------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Columns B-G hold numeric-looking values that are actually stored as
# plain text in this workbook (e.g. "3.10" keeps its trailing zero).
# Force Text format on those columns before assigning so Excel doesn't
# silently coerce them to numbers, then restore the default (no special
# format) style by pasting formats from an untouched cell.
$textRange = $q4.Range("B2:G6")
$textRange.NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "161724"
$q4.Range("C2").Value = "招商中证煤炭等权指数（LOF）A"
$q4.Range("D2").Value = "17.24"
$q4.Range("E2").Value = "93.84"
$q4.Range("F2").Value = "3.10"
$q4.Range("G2").Value = "0.5344"
$q4.Range("H2").Value = 6

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "013596"
$q4.Range("C3").Value = "招商中证煤炭等权指数（LOF）C"
$q4.Range("D3").Value = "1.56"
$q4.Range("E3").Value = "93.84"
$q4.Range("F3").Value = "3.10"
$q4.Range("G3").Value = "0.0484"
$q4.Range("H3").Value = 6

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "008736"
$q4.Range("C4").Value = "南方高股息主题股票A"
$q4.Range("D4").Value = "0.79"
$q4.Range("E4").Value = "91.27"
$q4.Range("F4").Value = "3.23"
$q4.Range("G4").Value = "0.0255"
$q4.Range("H4").Value = 7

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "016347"
$q4.Range("C5").Value = "招商中证煤炭等权指数（LOF）E"
$q4.Range("D5").Value = "0.20"
$q4.Range("E5").Value = "93.84"
$q4.Range("F5").Value = "3.10"
$q4.Range("G5").Value = "0.0062"
$q4.Range("H5").Value = 6

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "008737"
$q4.Range("C6").Value = "南方高股息主题股票C"
$q4.Range("D6").Value = "0.07"
$q4.Range("E6").Value = "91.27"
$q4.Range("F6").Value = "3.23"
$q4.Range("G6").Value = "0.0023"
$q4.Range("H6").Value = 7

# Restore default style (no NumberFormat override) on the text columns by
# pasting formats from a pristine, untouched cell on the same sheet.
$q4.Range("Z1:Z5").Copy()
$q4.Range("B2:B6").PasteSpecial(-4122)
$q4.Range("Z1:Z5").Copy()
$q4.Range("C2:C6").PasteSpecial(-4122)
$q4.Range("Z1:Z5").Copy()
$q4.Range("D2:D6").PasteSpecial(-4122)
$q4.Range("Z1:Z5").Copy()
$q4.Range("E2:E6").PasteSpecial(-4122)
$q4.Range("Z1:Z5").Copy()
$q4.Range("F2:F6").PasteSpecial(-4122)
$q4.Range("Z1:Z5").Copy()
$q4.Range("G2:G6").PasteSpecial(-4122)

$q4.Range("A1").Select()
